$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 8.182051697463905
    3 = 8.01007379356743
    4 = 8.797239129092313
    5 = 16.27321600439371
    6 = 8.066414953838478
    7 = 13.847952712232475
    8 = 8.066414953838478
    9 = 9.53112352006799
    10 = 24.439132215938336
    11 = 9.375665432811145
    12 = 9.509517250266017
    13 = 7.947518262177969
    14 = 1.6956331180355235
    15 = 26.72690099967158
    16 = 26.51960352917333
    17 = 1.888714214409287
    18 = 2.1161179869524958
    19 = 2.388644328245281
    20 = 2.532934269248223
    21 = 19.02960499761069
    22 = 10.320341398529736
    23 = 18.820685716437655
    24 = 0.02506131343921325
    25 = 1.1374866008016582
    26 = 2.4715826355868273
    27 = 19.665250046680914
    28 = 7.830657864141614
    29 = 10.025872472929654
    30 = 20.72044788298934
    31 = 7.314630481657307
    32 = 16.286917622833272
    33 = 21.475218246340425
    34 = 8.243427696238571
    35 = 8.373913822829692
    36 = 8.106373956737517
    37 = 8.055778489858271
    38 = 8.72494440346205
    39 = 8.116633350158384
    40 = 7.306611216007376
    41 = 18.484776081815134
    42 = 14.592042015506177
    43 = 9.408867436975976
    44 = 10.193821861256627
    45 = 10.407324039797578
    46 = 9.58022861008955
    47 = 9.2356417145884
    48 = 8.82987739958533
    49 = 10.536044417889133
    50 = 8.061966418738706
    51 = 20.182618982373334
    52 = 8.53929672348277
    53 = 8.045665895994912
    54 = 8.712526974050162
    55 = 7.143281865265328
    56 = 7.9735520826985145
    57 = 8.116408979950942
    58 = 8.126672533600978
    59 = 8.205673189421052
    60 = 8.113518481252118
    61 = 7.259661300889963
    62 = 5.665243484611729
    63 = 9.796931009270766
    64 = 8.55116233533741
    65 = 8.22340459814114
    66 = 8.19767477605039
    67 = 8.043622189151318
    68 = 8.041585393652957
    69 = 9.655701189819537
    70 = 8.124654988460366
    71 = 8.038089748485742
    72 = 9.295191874139414
    73 = 8.126672533600978
    74 = 8.279503628577823
    75 = 8.188415164548694
    76 = 8.064966541540802
    77 = 8.126672533600978
    78 = 8.041585393652957
    79 = 8.112797248765638
    80 = 9.158211014309318
    81 = 8.038089748485742
    82 = 1.0771238161772523
    83 = 1.4277459836095865
    84 = 26.764185565943876
    85 = 9.051311878931172
    86 = 7.960264118588261
    87 = 3.844910428363093
    88 = 2.8733073656649424
    89 = 5.382038655670804
    90 = 7.221737408815704
    91 = 8.085485626565141
    92 = 7.833028805737936
    93 = 3.289320882033337
    94 = 10.369883126128844
    95 = 8.122221111596561
    96 = 2.489495835460122
    97 = 1.3219314005900458
    98 = 20.241808453629286
    99 = 22.6937769504415
    100 = 8.712526974050162
    101 = 8.04014545100684
    102 = 8.007724151425661
    103 = 7.770721672306273
    104 = 1.6398439401583076
    105 = 28.36344187304499
    106 = 7.558956891721355
    107 = 10.11188311612278
    108 = 8.414682326070128
    109 = 7.522058747309277
    110 = 20.68283916500375
    111 = 7.274747862484206
    112 = 7.491042587564012
    113 = 20.303155578657908
    114 = 8.499619843263137
    115 = 8.755481845386008
    116 = 20.213788819687792
    117 = 8.058327007554423
    118 = 9.180512015107126
    119 = 8.261527884331175
    120 = 19.902527532373288
    121 = 25.947419935982627
    122 = 12.56686531163366
    123 = 17.82183397456314
    124 = 21.542950261359646
    125 = 21.797556850170416
    126 = -5.322175589390998
    127 = 19.325920804401452
    128 = 19.928636997784185
    129 = 19.706757154857623
    130 = 10.766543814641956
    131 = 8.672087666962952
    132 = 20.519673649038754
    133 = 6.238616055983375
    134 = 8.734476878105456
    135 = 10.600054037404046
    136 = 19.22901239798358
    137 = 19.7936017952748
    138 = 10.994586801409023
    139 = 8.95520482610032
    140 = 9.476101435602839
    141 = 8.480656440217233
    142 = 8.1101082688071
    143 = 8.042097232069027
    144 = 9.45881645135891
    145 = 9.559025944959814
    146 = 9.57971468702644
    147 = 10.5759652993475
    148 = 9.58985595465962
    149 = 9.015425045986094
    150 = 7.434115509907118
    151 = 8.19163028326565
    152 = 8.073739774328525
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $values[$row]
}
